$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cell text updates from the crypto price refresh.
# NumberFormat is forced to Text ("@") before assignment so that numeric-looking
# strings (e.g. "1.00", "0.640") are preserved verbatim as text instead of being
# coerced into numbers (which would drop trailing zeros / formatting).
# Style is reset back to "Normal" afterwards so no stray formatting is introduced.
function Set-TextValue($cellRef, $value) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue 'D2' '65.519.37'
Set-TextValue 'E2' '  +0.00%  '
Set-TextValue 'D3' '3.442.77'
Set-TextValue 'E3' '  +2.23%  '
Set-TextValue 'D4' '1.00'
Set-TextValue 'E4' '  +0.02%  '
Set-TextValue 'D5' '552.76'
Set-TextValue 'E5' '  +2.75%  '
Set-TextValue 'E6' '  -1.69%  '
Set-TextValue 'D7' '0.640'
Set-TextValue 'E7' '  +6.03%  '
Set-TextValue 'E8' '  -0.08%  '
Set-TextValue 'D9' '0.634'
Set-TextValue 'E9' '  +1.23%  '
Set-TextValue 'D10' '0.153'
Set-TextValue 'E10' '  +6.38%  '
Set-TextValue 'D11' '53.90'
Set-TextValue 'E11' '  -1.93%  '
Set-TextValue 'D12' '0.0000272'
Set-TextValue 'E12' '  +2.31%  '
Set-TextValue 'D13' '9.23'
Set-TextValue 'E13' '  +0.11%  '
Set-TextValue 'D14' '3.981.19'
Set-TextValue 'E14' '  +1.87%  '
Set-TextValue 'B15' 'Chainlink'
Set-TextValue 'C15' 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue 'D15' '18.46'
Set-TextValue 'E15' '  +2.68%  '
Set-TextValue 'B16' 'TRON'
Set-TextValue 'C16' 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue 'D16' '0.121'
Set-TextValue 'E16' '  +0.45%  '
Set-TextValue 'D17' '3.438.03'
Set-TextValue 'E17' '  +2.17%  '
Set-TextValue 'D18' '65.527.42'
Set-TextValue 'E18' '  -0.18%  '
Set-TextValue 'D19' '11.88'
Set-TextValue 'E19' '  +3.53%  '
Set-TextValue 'D20' '0.988'
Set-TextValue 'E20' '  +0.76%  '
Set-TextValue 'D21' '415.52'
Set-TextValue 'E21' '  +5.99%  '
Set-TextValue 'D22' '4.04'
Set-TextValue 'E22' '  +6.51%  '
Set-TextValue 'D23' '85.92'
Set-TextValue 'E23' '  +3.31%  '
Set-TextValue 'D24' '4.13'
Set-TextValue 'E24' '  -2.56%  '
Set-TextValue 'D25' '10.83'
Set-TextValue 'E25' '  -7.75%  '
Set-TextValue 'D26' '2.87'
Set-TextValue 'E26' '  +1.12%  '
Set-TextValue 'D27' '12.48'
Set-TextValue 'E27' '  +7.95%  '
Set-TextValue 'E28' '  -1.64%  '
Set-TextValue 'D29' '9.05'
Set-TextValue 'E29' '  +7.08%  '
Set-TextValue 'D30' '29.94'
Set-TextValue 'E30' '  +1.54%  '
Set-TextValue 'D31' '6.54'
Set-TextValue 'E31' '  -3.20%  '
Set-TextValue 'D32' '606.66'
Set-TextValue 'E32' '  -8.52%  '
Set-TextValue 'D33' '11.72'
Set-TextValue 'E33' '  +2.43%  '
Set-TextValue 'E34' '  +0.88%  '
Set-TextValue 'D35' '59.12'
Set-TextValue 'E35' '  +2.06%  '
Set-TextValue 'D36' '1.00'
Set-TextValue 'E36' '  +0.06%  '
Set-TextValue 'B37' 'PEPE'
Set-TextValue 'C37' 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue 'D37' '0.0₃0794'
Set-TextValue 'E37' '  +2.35%  '
Set-TextValue 'D38' '37.51'
Set-TextValue 'E38' '  -0.18%  '
Set-TextValue 'B39' 'Kaspa'
Set-TextValue 'C39' 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue 'D39' '0.145'
Set-TextValue 'E39' '  +12.51%  '
Set-TextValue 'D40' '0.380'
Set-TextValue 'E40' '  -4.14%  '
Set-TextValue 'D41' '3.209.08'
Set-TextValue 'E41' '  +6.58%  '
Set-TextValue 'D42' '3.35'
Set-TextValue 'E42' '  +2.22%  '
Set-TextValue 'D43' '0.996'
Set-TextValue 'E43' '  -0.33%  '
Set-TextValue 'D44' '2.56'
Set-TextValue 'E44' '  -7.25%  '
Set-TextValue 'D45' '2.82'
Set-TextValue 'E45' '  +1.58%  '
Set-TextValue 'D46' '0.0413'
Set-TextValue 'E46' '  +0.49%  '
Set-TextValue 'D47' '3.24'
Set-TextValue 'E47' '  +0.22%  '
Set-TextValue 'D48' '2.71'
Set-TextValue 'E48' '  -0.60%  '
Set-TextValue 'D49' '0.133'
Set-TextValue 'E49' '  +3.95%  '
Set-TextValue 'D50' '138.07'
Set-TextValue 'E50' '  -1.33%  '
Set-TextValue 'D51' '8.41'
Set-TextValue 'E51' '  -1.49%  '
